# Apply the "gh-pages output generated at 456a3b4" update to
# 北京-漫展信息.xlsx (Beijing convention/expo info workbook).
#
# Sheets: 展览(1)=Exhibitions, 演出(2)=Performances,
#         本地生活(3)=Local life, 全部类型(4)=All types

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: 展览 (Exhibitions)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# "Want to go" count bumps that do NOT involve any row shifting
# (rows 1-30, before the newly inserted row).
$ws1.Range("F5").Value  = 357
$ws1.Range("F6").Value  = 408
$ws1.Range("F7").Value  = 907
$ws1.Range("F8").Value  = 69
$ws1.Range("F9").Value  = 556
$ws1.Range("F15").Value = 46
$ws1.Range("F17").Value = 6741
$ws1.Range("F18").Value = 68
$ws1.Range("F19").Value = 79
$ws1.Range("F20").Value = 25
$ws1.Range("F21").Value = 7667
$ws1.Range("F24").Value = 3425
$ws1.Range("F25").Value = 36
$ws1.Range("F26").Value = 2160
$ws1.Range("F29").Value = 194
$ws1.Range("F30").Value = 355

# A brand-new event was added, inserted as row 35, pushing the
# previous rows 35-45 down to 36-46.
$ws1.Range("A35").EntireRow.Insert()

$ws1.Range("A35").Value = 34
$ws1.Range("A35").Font.Bold = $true
$ws1.Range("A35").HorizontalAlignment = -4108
$ws1.Range("A35").VerticalAlignment = -4160
$ws1.Range("A35").Borders.LineStyle = 1

$ws1.Range("B35").Value = "2024.04.04"
$ws1.Range("C35").Value = "北京·IDOx梦次元动漫游戏嘉年华3rd同人创作大会"
$ws1.Range("D35").Value = "北京展览馆 北京展览馆"
$ws1.Range("E35").Value = "2024.04.04 09:30-04.05 17:00"
$ws1.Range("F35").Value = 0
$ws1.Range("G35").Value = "不可售"
$ws1.Range("H35").Value = "https://show.bilibili.com/platform/detail.html?id=82023"
$ws1.Range("I35").Value = "//i2.hdslb.com/bfs/openplatform/202402/DE1Xw4Ne1708668500346.png"

# Rows 36-46 now hold what used to be rows 35-45; most of the text
# carried straight across with the insert, only the "want to go"
# counts (and one title rename) need touching up.
$ws1.Range("F36").Value = 1792   # 第15届IJOY漫展xCGF游戏节

# row 37 (XW无限世界cosplay嘉年华) unchanged

$ws1.Range("F38").Value = 204    # Yok运动番Only
$ws1.Range("F39").Value = 59     # thebONE GOJO超次元嘉年华12nd

$ws1.Range("C40").Value = "北京·QMQ动漫游戏嘉年华"
$ws1.Range("F40").Value = 11

# row 41 (次元风暴游园会) unchanged

$ws1.Range("F42").Value = 1261   # IDO动漫游戏嘉年华45th
$ws1.Range("F43").Value = 4      # IDO动漫游戏嘉年华45th同人创作大会
$ws1.Range("F44").Value = 1907   # 第16届IJOY漫展XCGF游戏节

# row 45 (原神only3.0) unchanged
# row 46 (次元风暴游园会2.0) unchanged

# ---------------------------------------------------------------
# Sheet 2: 演出 (Performances)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 21
$ws2.Range("F7").Value = 84

# ---------------------------------------------------------------
# Sheet 3: 本地生活 (Local life)
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("E2").Value = "2023.12.22 10:00-2024.03.15 17:00"

# ---------------------------------------------------------------
# Sheet 4: 全部类型 (All types)
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("E2").Value = "2023.12.22 10:00-2024.03.15 17:00"

$ws4.Range("F7").Value  = 357
$ws4.Range("F8").Value  = 408
$ws4.Range("F9").Value  = 907
$ws4.Range("F10").Value = 69
$ws4.Range("F11").Value = 556
$ws4.Range("F18").Value = 46
$ws4.Range("F20").Value = 6741
$ws4.Range("F21").Value = 68
$ws4.Range("F22").Value = 79
$ws4.Range("F23").Value = 25
$ws4.Range("F24").Value = 7667
$ws4.Range("F27").Value = 3425
$ws4.Range("F28").Value = 36
$ws4.Range("F29").Value = 2160
$ws4.Range("F32").Value = 194
$ws4.Range("F33").Value = 355
$ws4.Range("F38").Value = 1792

$ws4.Range("F40").Value = 204
$ws4.Range("F41").Value = 59

$ws4.Range("C42").Value = "北京·QMQ动漫游戏嘉年华"
$ws4.Range("F42").Value = 11

$ws4.Range("F44").Value = 1261
$ws4.Range("F45").Value = 1907
$ws4.Range("F46").Value = 21
$ws4.Range("F49").Value = 84

Write-Host "Edit complete"
